$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.564.21"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "2.677.11"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "598.10"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "166.61"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.26%  "

$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").Value = "2.675.33"
$ws.Range("E9").Value = "  -0.67%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.143"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.43%  "

$ws.Range("E11").Value = "  +1.28%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.359"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("E13").Value = "  -1.40%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.83"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.47%  "

$ws.Range("D15").Value = "3.163.77"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "67.498.62"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "2.672.00"
$ws.Range("E18").Value = "  -1.04%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.74"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "363.75"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.38"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.38%  "

$ws.Range("E23").Value = "  -1.38%  "

$ws.Range("E24").Value = "  -4.04%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -4.70%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("E30").Value = "  -0.10%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "557.85"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.03"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("E33").Value = "  -3.42%  "

$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -4.92%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "19.54"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.35%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "154.78"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("E42").Value = "  -4.07%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "17.95"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("E44").Value = "  +0.02%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.52"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.36%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "40.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.11%  "

$ws.Range("D47").Value = "0.0₆0298"
$ws.Range("E47").Value = "  -5.63%  "

$ws.Range("E48").Value = "  -1.69%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "153.36"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("E50").Value = "  -2.21%  "

